$wb = $excel.ActiveWorkbook

# --- Add the new "Scope" worksheet as the last tab (after "Tables") ---
$scope = $wb.Worksheets.Add()
$scope.Name = "Scope"
$scope.Move($wb.Worksheets.Count)
$scope = $wb.Worksheets.Item("Scope")

# --- Populate "Scope" sheet content ---
$scope.Range("B2").Value = "Phase 1"

$scope.Range("C4").Value = "Develop application lookup maintenance forms:"
$scope.Range("J4").Value = "Other tables"
$scope.Range("D5").Value = "User maintenance"
$scope.Range("D6").Value = "User Level maintenance"
$scope.Range("J6").Value = "GST table - for now, embedded in Customer table"
$scope.Range("J6").Interior.Color = 65535
$scope.Range("D7").Value = "Order Status maintenance"
$scope.Range("D8").Value = "Product Category maintenance"

$scope.Range("C10").Value = "Develop application main forms:"
$scope.Range("D11").Value = "Customer maintenance"
$scope.Range("D12").Value = "Product maintenance"
$scope.Range("D13").Value = "Supplier maintenance"
$scope.Range("D14").Value = "Order entry form"

$scope.Range("B17").Value = "Phase 1"

$scope.Range("C19").Value = "Develop "
$scope.Range("D20").Value = "Generate order paper/pdf form"
$scope.Range("D21").Value = "Generate order fulfillment checklist"
$scope.Range("D22").Value = "Generate invoice paper "
$scope.Range("D23").Value = "Login form"
$scope.Range("D24").Value = "Order processing form"
$scope.Range("H24").Value = "This form is for order fulfillment, shipping order"
$scope.Range("E25").Value = "Order processing backend"
$scope.Range("H25").Value = "This will be responsible for updating inventory based on order once fulfilled"
$scope.Range("D26").Value = "Inventory maintenance"
$scope.Range("H26").Value = "This form is for updating product inventory as  products are received for top up.  Order count adjustment process, etc."

$scope.Range("B31").Value = "Phase 3"
$scope.Range("C32").Value = "Develop functionality to send pdf order form to customer via email, from remote or backoffice"
$scope.Range("C33").Value = "Develop remote order sync/upload to main database process"
$scope.Range("J33").Value = "This might require additional columns in the main database, orders table, for keeping track of status of sync/upload"
$scope.Range("C34").Value = "Develop data sync from main database to remote units"
$scope.Range("J34").Value = "This process is for refreshing data in remote units.  Example, new customer created by admin in main database."
$scope.Range("P35").Value = "new products created by admin in main database."
$scope.Range("P36").Value = "and possibly other type of data."

# --- Update the "Overview" sheet view: clear scroll position, zoom to 80%, change selection ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Activate()
$excel.ActiveWindow.Zoom = 80
$overview.Range("F34").Select()

# --- Make "Scope" the active/selected tab, matching the saved workbook state ---
$scope.Activate()
$scope.Range("T24").Select()
